# Adds the new tank-experiment log rows (74-93) to Sayfa1, matching the
# "new excel sheets added to improve performance" data-entry commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$newRows = @(
    @{Row=74; A='40 lpm - hole 6 + 7 - diameter 12 (Boğaza Bağlı)'; B='17:12:41'; C='23.08.2021'},
    @{Row=75; A='40 lpm - hole 4 + 5 - diameter 12 (Boğaza Bağlı)'; B='08:41:49'; C='24.08.2021'},
    @{Row=76; A='40 lpm - hole 5 + 7 - diameter 12 (Boğaza Bağlı)'; B='08:57:36'; C='24.08.2021'},
    @{Row=77; A='40 lpm - hole 5 + 6 - diameter 12 (Boğaza Bağlı)'; B='09:16:58'; C='24.08.2021'},
    @{Row=78; A='40 lpm - hole 3 + 7 - diameter 12 (Boğaza Bağlı)'; B='09:38:16'; C='24.08.2021'},
    @{Row=79; A='80 lpm - hole 6 + 7 - diameter 12 (Boğaza Bağlı)'; B='08:27:03'; C='25.08.2021'},
    @{Row=80; A='80 lpm - hole 4 + 5 - diameter 12 (Boğaza Bağlı)'; B='08:42:45'; C='25.08.2021'},
    @{Row=81; A='80 lpm - hole 5 + 7 - diameter 12 (Boğaza Bağlı)'; B='09:05:28'; C='25.08.2021'},
    @{Row=82; A='80 lpm - hole 5 + 6 - diameter 12 (Boğaza Bağlı)'; B='07:49:04'; C='27.08.2021'},
    @{Row=83; A='80 lpm - hole 3 + 7 - diameter 12 (Boğaza Bağlı)'; B='08:05:59'; C='27.08.2021'},
    @{Row=84; A='120 lpm - hole 6 + 7 - diameter 12 (Boğaza Bağlı)'; B='09:51:15'; C='27.08.2021'},
    @{Row=85; A='120 lpm - hole 4 + 5 - diameter 12 (Boğaza Bağlı)'; B='10:29:22'; C='27.08.2021'},
    @{Row=86; A='120 lpm - hole 5 + 7 - diameter 12 (Boğaza Bağlı)'; B='10:48:07'; C='27.08.2021'},
    @{Row=87; A='120 lpm - hole 5 + 6 - diameter 12 (Boğaza Bağlı)'; B='11:07:14'; C='27.08.2021'},
    @{Row=88; A='120 lpm - hole 3 + 7 - diameter 12 (Boğaza Bağlı)'; B='11:36:18'; C='27.08.2021'},
    @{Row=89; A='120 lpm - hole 1 + 2 - diameter 12 (Boğaza Bağlı)'; B='14:07:43'; C='27.08.2021'},
    @{Row=90; A='120 lpm - hole 1 + 3 - diameter 12 (Boğaza Bağlı)'; B='14:24:27'; C='27.08.2021'},
    @{Row=91; A='120 lpm - hole 1 + 4 - diameter 12 (Boğaza Bağlı)'; B='14:45:06'; C='27.08.2021'},
    @{Row=92; A='120 lpm - hole 2 + 3 - diameter 12 (Boğaza Bağlı)'; B='16:08:54'; C='27.08.2021'},
    @{Row=93; A='120 lpm - hole 2 + 4 - diameter 12 (Boğaza Bağlı)'; B='16:25:49'; C='27.08.2021'}
)

foreach ($r in $newRows) {
    # Keep the new B/C cells unstyled (no explicit xf) just like every other
    # data row on the sheet - they simply inherit the column's left-align
    # style rather than carrying their own style index.
    $ws.Cells.Item($r.Row, 2).Style = "Normal"
    $ws.Cells.Item($r.Row, 3).Style = "Normal"

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

# Widen column A to fit the longer "(Boğaza Bağlı)" descriptions.
$ws.Columns.Item(1).ColumnWidth = 41.3

# Scroll/zoom/selection state the author left the sheet in after the edit.
$ws.Application.ActiveWindow.Zoom = 99
$ws.Range("F84").Select()
